$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-style the existing "Temps joué" (G) and stat (H:V) columns for rows 307-320 ---
# G column: left-align, vertical-center
$ws.Range("G307:G320").HorizontalAlignment = -4131
$ws.Range("G307:G320").VerticalAlignment = -4108
# H:V columns: right-align, vertical-center
$ws.Range("H307:V320").HorizontalAlignment = -4152
$ws.Range("H307:V320").VerticalAlignment = -4108

# --- Append the six new "Entrainement" rows (321-326) for the J-1 session ---
# Give column B (Date) the same date-serial formatting as the row above it before
# writing the new date values, so the new cells pick up the existing date style.
$ws.Range("B320").Copy($ws.Range("B321:B326"))

$newRows = @(
    @{ Row=321; A="Entrainement"; B=45892; C="Global"; D="J-1"; E="Romain Thunet";      F="center back";      G="01:21:22"; H=4.63; I=0.1;  J=4.53;               K=0.1;  L=0;                    M=0;    N=0; O=0;  P=3.3;  Q=19.34; R=4.45;               S=29; T=3;  U=22; V=2  },
    @{ Row=322; A="Entrainement"; B=45892; C="Global"; D="J-1"; E="Karahali Souaré";    F="right forward";    G="01:19:17"; H=5.66; I=0.32; J=5.33;               K=0.25; L=0.07;                  M=0.01; N=0; O=2;  P=3.84; Q=26.45; R=5.27;               S=39; T=10; U=37; V=19 },
    @{ Row=323; A="Entrainement"; B=45892; C="Global"; D="J-1"; E="Ilyes Boughanmi";    F="center forward";   G="01:19:48"; H=5.01; I=0.23; J=4.7699999999999996; K=0.22; L=0.02;                  M=0;    N=0; O=0;  P=3.68; Q=21.14; R=5.4;                S=29; T=12; U=16; V=5  },
    @{ Row=324; A="Entrainement"; B=45892; C="Global"; D="J-1"; E="Omar Benyounes";     F="center midfield";  G="01:21:38"; H=5.48; I=0.36; J=5.1100000000000003; K=0.3;  L=0.06;                  M=0.01; N=0; O=1;  P=3.93; Q=25.97; R=4.67;               S=30; T=6;  U=29; V=5  },
    @{ Row=325; A="Entrainement"; B=45892; C="Global"; D="J-1"; E="Amine Taiar";        F="center back";      G="01:20:20"; H=5.09; I=0.21; J=4.87;               K=0.17; L=0.04;                  M=0;    N=0; O=0;  P=3.69; Q=24.15; R=4.53;               S=22; T=3;  U=21; V=5  },
    @{ Row=326; A="Entrainement"; B=45892; C="Global"; D="J-1"; E="Ilan Ihaddadene";    F="center midfield";  G="01:19:48"; H=5.64; I=0.3;  J=5.33;               K=0.21; L=0.1;                   M=0;    N=0; O=0;  P=4.1399999999999997; Q=24.55; R=4.6399999999999997; S=30; T=5;  U=23; V=6  }
)

foreach ($rd in $newRows) {
    $r = $rd.Row
    $ws.Cells.Item($r, 1).Value = $rd.A
    $ws.Cells.Item($r, 2).Value = $rd.B
    $ws.Cells.Item($r, 3).Value = $rd.C
    $ws.Cells.Item($r, 4).Value = $rd.D
    $ws.Cells.Item($r, 5).Value = $rd.E
    $ws.Cells.Item($r, 6).Value = $rd.F
    $ws.Cells.Item($r, 7).Value = $rd.G
    $ws.Cells.Item($r, 8).Value = $rd.H
    $ws.Cells.Item($r, 9).Value = $rd.I
    $ws.Cells.Item($r, 10).Value = $rd.J
    $ws.Cells.Item($r, 11).Value = $rd.K
    $ws.Cells.Item($r, 12).Value = $rd.L
    $ws.Cells.Item($r, 13).Value = $rd.M
    $ws.Cells.Item($r, 14).Value = $rd.N
    $ws.Cells.Item($r, 15).Value = $rd.O
    $ws.Cells.Item($r, 16).Value = $rd.P
    $ws.Cells.Item($r, 17).Value = $rd.Q
    $ws.Cells.Item($r, 18).Value = $rd.R
    $ws.Cells.Item($r, 19).Value = $rd.S
    $ws.Cells.Item($r, 20).Value = $rd.T
    $ws.Cells.Item($r, 21).Value = $rd.U
    $ws.Cells.Item($r, 22).Value = $rd.V
}

# --- Match the author's final viewport/selection state ---
$ws.Range("E330").Select()
